# Applies the "add new SCH Sheet" edit to 核心器件.docx
#
# The change set (per the OOXML diff) touches five paragraphs:
#   1. "振动路径前置放大：AD797(SOP8)"      -> drop the paragraph-mark rFonts hint (pPr removed)
#   2. "单端转差分：ADA4041-1(SOP8)"        -> "ADA4041-1(SOP8)" splits into "ADA49" + "41-1(SOP8)"
#   3. "485通信接口芯片：MAX3461CSD（SOP14）" -> "MAX3461"/"CSD"/"SOP14" become "MAX"/"485"/"SOP8"
#   4. "电源+5V转3.3V：LM1117-3"            -> appends a new run ".3"
#   5. "双排母，18p，10p"                   -> drop the paragraph-mark rFonts hint (pPr removed)
#                                             and insert "2X" runs before "18p" and before "10p"
#
# Because several of these edits split/merge runs and one even removes the
# paragraph-mark run-properties (<w:pPr><w:rPr>...), plain Find/Replace can't
# reproduce the exact run layout (Find's replace coalesces the whole
# paragraph into one run). Instead each affected paragraph's whole Range
# (which includes its end-of-paragraph mark) is swapped out in one shot via
# Range.InsertXML with a minimal WordprocessingML package fragment that
# spells out precisely the runs/pPr we want - InsertXML replaces the full
# paragraph when the target Range spans it end-to-end, which is exactly the
# granularity we need here.

$d = $word.ActiveDocument

function New-PackageXml([string]$bodyXml) {
    $header = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
    $footer = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    return $header + $bodyXml + $footer
}

function Replace-ParagraphByMarker([string]$marker, [string]$newParaXml) {
    $target = $null
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text -like "*$marker*") {
            $target = $p
            break
        }
    }
    if ($target -eq $null) {
        Write-Host "WARNING: paragraph matching '$marker' not found"
        return
    }
    $pkg = New-PackageXml $newParaXml
    $target.Range.InsertXML($pkg)
}

# 1) "振动路径前置放大：AD797(SOP8)" - remove the paragraph-mark rFonts hint (pPr)
$p1New = '<w:p w:rsidR="004358AB" w:rsidRDefault="00EA66D3" w:rsidP="00C40D41"><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>振动路径</w:t></w:r><w:r w:rsidR="008603F9"><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>前置</w:t></w:r><w:r w:rsidR="00C40D41"><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>放大：</w:t></w:r><w:r w:rsidR="00C40D41"><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>AD797</w:t></w:r><w:r w:rsidR="009B57B7"><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>(SOP8)</w:t></w:r></w:p>'
Replace-ParagraphByMarker "振动路径前置放大" $p1New

# 2) "单端转差分：ADA4041-1(SOP8)" - split the last run into "ADA49" + "41-1(SOP8)"
$p2New = '<w:p w:rsidR="00873262" w:rsidRDefault="00873262" w:rsidP="00C40D41"><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>单端转差分：</w:t></w:r><w:r w:rsidR="007D680E"><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>ADA49</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>41-1(SOP8)</w:t></w:r></w:p>'
Replace-ParagraphByMarker "单端转差分" $p2New

# 3) "485通信接口芯片：MAX3461CSD（SOP14）" -> "MAX3461"/"CSD"/"SOP14" become "MAX"/"485"/"SOP8"
$p3New = '<w:p w:rsidR="00AD44D9" w:rsidRDefault="002861B7" w:rsidP="00C40D41"><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>485</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>通信接口芯片：</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>MAX</w:t></w:r><w:r w:rsidR="005743FE"><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>485</w:t></w:r><w:r w:rsidR="00D65D44"><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>（</w:t></w:r><w:r w:rsidR="00D65D44"><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>SOP8</w:t></w:r><w:r w:rsidR="00D65D44"><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>）</w:t></w:r></w:p>'
Replace-ParagraphByMarker "通信接口芯片：MAX" $p3New

# 4) "电源+5V转3.3V：LM1117-3" - append a new run ".3"
$p4New = '<w:p w:rsidR="00C967E4" w:rsidRDefault="006F280C" w:rsidP="00322D81"><w:pPr><w:ind w:left="720" w:hanging="720"/></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>电源</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>+5V</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>转</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>3.3V</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>：</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>LM1117-3</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>.3</w:t></w:r></w:p>'
Replace-ParagraphByMarker "3.3V" $p4New

# 5) "双排母，18p，10p" - drop the paragraph-mark rFonts hint (pPr) and insert "2X" before each count
$p5New = '<w:p w:rsidR="00B03BC6" w:rsidRDefault="0032537C" w:rsidP="00750FBA"><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>双排母，</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>2X</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>18p</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>，</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>2X</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>10p</w:t></w:r></w:p>'
Replace-ParagraphByMarker "双排母" $p5New

Write-Host "Done."
